# WPCT_Fig_7.6_RequirementsSpecification.docx
#
# Applies the "Fig. 7.7 implemented in SAS" commit's minor, for-consistency
# edits to the Fig. 7.6 spec: wherever the spec singled out "Last Baseline"
# / "Last Post-baseline", it now also allows "Minimum"/"Maximum" (i.e.
# "Last/Min/Max ..."), and the box-plot grouping variable changes from
# AVISITN to STUDYID.

$d = $word.ActiveDocument

function Replace-ViaRange($Find, $Replace) {
    $rng = $d.Content.Duplicate
    $ok = $rng.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $Find"
        return $false
    }
    # Assign the literal replacement text (rather than using Find's own
    # Replace mode) so straight quotes in the replacement aren't mangled
    # into curly "smart" quotes.
    $rng.Text = $Replace
    return $true
}

# 1. Title: "... at Last Baseline and Last Post-baseline for Multiple
#    Studies and Analysis Timepoint" -> "... at Last/Min/Max Baseline and
#    Last/Min/Max Post-baseline for Multiple Studies and Analysis Timepoint"
Replace-ViaRange "at Last Baseline and Last Post-baseline for Multiple Studies and Analysis Timepoint" `
                  "at Last/Min/Max Baseline and Last/Min/Max Post-baseline for Multiple Studies and Analysis Timepoint"

# 2. Bullet: "Display results for just two study visits: Last Baseline and
#    Last Post-baseline" -> "... Last/Min/Max Baseline and Last/Min/Max
#    Post-baseline"
Replace-ViaRange "Display results for just two study visits: Last Baseline and Last Post-baseline" `
                  "Display results for just two study visits: Last/Min/Max Baseline and Last/Min/Max Post-baseline"

# 3. Bullet: "Clearly label and visually separate Last Baseline and Last
#    Post-baseline results" -> "Clearly label and separate Last/Min/Max
#    Baseline and Last/Min/Max Post-baseline results"
Replace-ViaRange "Clearly label and visually separate Last Baseline and Last Post-baseline results" `
                  "Clearly label and separate Last/Min/Max Baseline and Last/Min/Max Post-baseline results"

# 4. "Generates a box plot of AVAL by AVISITN and TRTPN." -> "... by
#    STUDYID and TRTPN." This is the site of the document's most recent
#    real edit, so Word's "_GoBack" bookmark (it tracks the last edit
#    location, and only ever exists once per document) needs to move here
#    from its old spot (in "HI, AVISITN" further down).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
Replace-ViaRange "Generates a box plot of AVAL by AVISITN and TRTPN." `
                  "Generates a box plot of AVAL by STUDYID and TRTPN."
$goBackRng = $d.Content.Duplicate
$goBackRng.Find.Execute("STUDYID", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRng.Start = $goBackRng.End
$d.Bookmarks.Add("_GoBack", $goBackRng)

# 5. '"Baseline" visit number for "Last Baseline"' -> '"Baseline" visit
#    number for "Last", "Minimum" or "Maximum" Baseline'
Replace-ViaRange '"Baseline" visit number for "Last Baseline"' `
                  '"Baseline" visit number for "Last", "Minimum" or "Maximum" Baseline'

# 6. '"Endpoint" visit number for "Last Post-baseline"' -> '"Endpoint"
#    visit number for "Last", "Minimum" or "Maximum" Post-baseline'
Replace-ViaRange '"Endpoint" visit number for "Last Post-baseline"' `
                  '"Endpoint" visit number for "Last", "Minimum" or "Maximum" Post-baseline'

# 7. Footer page-number field: cached PAGE display "3" -> "1" (document now
#    stands alone, so it is page 1 of its own footer pagination). Leave the
#    NUMPAGES field's cached "3" alone -- editing the PAGE field's cached
#    result via Find/Replace also refreshes NUMPAGES's cache in this
#    engine, so restore it afterwards.
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Find.Execute("3", $true, $false, $false, $false, $false, $true, 0, $false, "1", 2) | Out-Null
$numPagesField = $d.Sections.Item(1).Footers.Item(1).Range.Fields.Item(2)
$numPagesField.Result.Find.Execute("1", $true, $false, $false, $false, $false, $true, 0, $false, "3", 2) | Out-Null

Write-Output "Done"
